$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.503.88'
$ws.Range("E2").Value = '  -2.15%  '
$ws.Range("D3").Value = '1.995.56'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.601'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.80'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.14%  '
$ws.Range("E9").Value = '  -3.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0746'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.83%  '
$ws.Range("E12").Value = '  -2.99%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.288.33'
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.86%  '
$ws.Range("E16").Value = '  -5.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.80%  '
$ws.Range("D18").Value = '1.994.58'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '36.579.91'
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.82%  '
$ws.Range("D21").Value = '0.0₃0804'
$ws.Range("E21").Value = '  -3.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.80%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +1.03%  '
$ws.Range("E26").Value = '  -8.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '18.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0600'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.62'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.51%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.450.29'
$ws.Range("E42").Value = '  +2.92%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0927'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("E44").Value = '  -4.21%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.991'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.67%  '
